$d = $word.ActiveDocument

function Split-RunAt($pos) {
    # Forces a run boundary immediately before the absolute character
    # position $pos by toggling a character-level format (Bold on, then
    # back off) on the single character that precedes $pos. The engine
    # materializes a distinct run for that character without merging it
    # back into its neighbour, while the run keeps its original (computed)
    # formatting once Bold is restored to its original value.
    $r1 = $d.Range($pos - 1, $pos)
    $r1.Bold = 1
    $r2 = $d.Range($pos - 1, $pos)
    $r2.Bold = 0
}

# ---------------------------------------------------------------------
# 1) ".1 Diagrammes des cas d'utilisation" -> split off the final "n"
#    into its own run (text itself is unchanged).
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(37)
$end1 = $p1.Range.End
Split-RunAt($end1 - 1)

# ---------------------------------------------------------------------
# 2) ".2 Description détaillée" -> ".2 Diagrammes de séquences système"
#    split into three runs: "." | "2" | " Diagrammes de séquences système"
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(38)
$p2.Range.Find.Execute("Description détaillée", $false, $false, $false, $false, $false, $true, 1, $false, "Diagrammes de séquences système", 2) | Out-Null

$p2b = $d.Paragraphs.Item(38)
$start2 = $p2b.Range.Start
Split-RunAt($start2 + 8)
Split-RunAt($start2 + 9)

# ---------------------------------------------------------------------
# 3) ".3 Diagrammes de séquences système" -> ".3 Diagrammes de classe"
#    split into three runs: "." | "3" | " Diagrammes de classe"
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(39)
$p3.Range.Find.Execute("Diagrammes de séquences système", $false, $false, $false, $false, $false, $true, 1, $false, "Diagrammes de classe", 2) | Out-Null

$p3b = $d.Paragraphs.Item(39)
$start3 = $p3b.Range.Start
Split-RunAt($start3 + 8)
Split-RunAt($start3 + 9)

# ---------------------------------------------------------------------
# 4) Remove the old ".4 Diagrammes de classe" paragraph entirely
#    (its content has just been folded into paragraph 39 above), merging
#    the following paragraph ("2.2.2. Modélisation") up.
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(40)
$r4 = $d.Range($p4.Range.Start, $p4.Range.End)
$r4.Delete()
